$wb = $excel.ActiveWorkbook

# --- Update existing sheet "tc1" data (credentials changed) ---
$ws1 = $wb.Worksheets.Item("tc1")
$ws1.Range("A2").Value = "bhanu"
$ws1.Range("B2").Value = "akashara"

# --- Add new worksheet "ValidLogin" right after "tc1" ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "ValidLogin"

# Populate header row
$ws2.Range("A1").Value = "username"
$ws2.Range("B1").Value = "password"

# Populate data row - set B2 before A2 so shared-string insertion order
# matches (pointofsale, then ADMIN)
$ws2.Range("B2").Value = "pointofsale"
$ws2.Range("A2").Value = "ADMIN"

# --- View/selection state ---
# tc1 is no longer the selected tab; its selection moves to B3
[void]$ws1.Range("B3").Select()

# ValidLogin becomes the active/selected sheet with selection at A2
[void]$ws2.Activate()
[void]$ws2.Range("A2").Select()
$excel.ActiveWindow.Zoom = 205

Write-Host "Done"
